# MaterMinds TestCases.xlsx update
# - Add a new test case row (row 4): "Clicking on boxview color"
# - Resize columns A:E to fit the new content (best-fit-style column widths)
# - Move selection to the new last cell (E4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Clicking on boxview color"
$ws.Range("C4").Value = "should change the color depending on how many times the user clicks"
$ws.Range("D4").Value = "cycles throught colors each time user clicks "
$ws.Range("E4").Value = "PASS"

# Best-fit style column widths for the updated data
$ws.Columns.Item(1).ColumnWidth = 5.498697916666667
$ws.Columns.Item(2).ColumnWidth = 20.998697916666668
$ws.Columns.Item(3).ColumnWidth = 58.830729166666664
$ws.Columns.Item(4).ColumnWidth = 36.830729166666664
$ws.Columns.Item(5).ColumnWidth = 7.830729166666667

# Move the active selection onto the newly added row's last cell
$ws.Range("E4").Select()
